$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in C5 from 25 to 15
$ws.Range("C5").Value = 15

# Update the selection to just C5 (instead of B2:C5)
$ws.Range("C5").Select()
